$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "27÷6=4, 3"
$t.Cell(1, 2).Range.Text = "33÷7=4, 5"
$t.Cell(1, 3).Range.Text = "49÷7=7, 0"
$t.Cell(1, 4).Range.Text = "88÷7=12, 4"
$t.Cell(1, 5).Range.Text = "86÷4=21, 2"
$t.Cell(5, 1).Range.Text = "75÷3=25, 0"
$t.Cell(5, 2).Range.Text = "34÷4=8, 2"
$t.Cell(5, 3).Range.Text = "10÷4=2, 2"
$t.Cell(5, 4).Range.Text = "22÷3=7, 1"
$t.Cell(5, 5).Range.Text = "47÷6=7, 5"
$t.Cell(9, 1).Range.Text = "32÷6=5, 2"
$t.Cell(9, 2).Range.Text = "30÷2=15, 0"
$t.Cell(9, 3).Range.Text = "92÷2=46, 0"
$t.Cell(9, 4).Range.Text = "75÷9=8, 3"
$t.Cell(9, 5).Range.Text = "37÷8=4, 5"
$t.Cell(13, 1).Range.Text = "60÷5=12, 0"
$t.Cell(13, 2).Range.Text = "48÷9=5, 3"
$t.Cell(13, 3).Range.Text = "65÷2=32, 1"
$t.Cell(13, 4).Range.Text = "99÷7=14, 1"
$t.Cell(13, 5).Range.Text = "16÷7=2, 2"
$t.Cell(17, 1).Range.Text = "17÷8=2, 1"
$t.Cell(17, 2).Range.Text = "70÷9=7, 7"
$t.Cell(17, 3).Range.Text = "91÷2=45, 1"
$t.Cell(17, 4).Range.Text = "18÷9=2, 0"
$t.Cell(17, 5).Range.Text = "20÷4=5, 0"
